$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated GX/LX MSRP values
$ws.Range("D29").Value = 53100
$ws.Range("D30").Value = 55890
$ws.Range("D31").Value = 64365
$ws.Range("D32").Value = 86580
$ws.Range("E32").Value = 1025
$ws.Range("D33").Value = 91580
$ws.Range("E33").Value = 1025

# D34 previously held a blank placeholder string; it now gets a real MSRP
# value, so give it the same number format ($#,##0 style) as the other
# BASE MSRP cells in column D.
$ws.Range("D34").Value = 94475
$ws.Range("D34").NumberFormat = "#,##0"
$ws.Range("E34").Value = 1025
